# Rename all 50 worksheets to new summ<id> names (importance ranking / shap by postcode for distance scripts)
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ40041436",
    "summ40182768",
    "summ40363151",
    "summ40570101",
    "summ40806548",
    "summ41020364",
    "summ41205586",
    "summ41365224",
    "summ41591845",
    "summ41747631",
    "summ41889231",
    "summ42054152",
    "summ42212587",
    "summ42350770",
    "summ42526028",
    "summ42694506",
    "summ42850584",
    "summ43025532",
    "summ43178481",
    "summ43325528",
    "summ43465390",
    "summ43616675",
    "summ43768482",
    "summ43947239",
    "summ44127837",
    "summ44285168",
    "summ44439025",
    "summ44593681",
    "summ44742849",
    "summ44890256",
    "summ45054989",
    "summ45212341",
    "summ45373336",
    "summ45534587",
    "summ45694821",
    "summ45850824",
    "summ46009773",
    "summ46161770",
    "summ46311448",
    "summ46463169",
    "summ46648487",
    "summ46812726",
    "summ46989034",
    "summ47144097",
    "summ47308546",
    "summ47477103",
    "summ47638275",
    "summ47800076",
    "summ47957559",
    "summ48138557"
)

for ($i = 1; $i -le $newNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

